# Apply crypto price/volume updates per commit: "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.738.65"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").Value = "2.528.45"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'309.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").Value = "'100.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("E7").Value = "  -1.10%  "
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("D10").Value = "'35.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").Value = "'0.0805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.02%  "
$ws.Range("D12").Value = "'7.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "2.917.31"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.554.57"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'15.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.88%  "
$ws.Range("D17").Value = "'0.813"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.04%  "
$ws.Range("D18").Value = "42.721.20"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "'6.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.34%  "
$ws.Range("D20").Value = "0.0₃0950"
$ws.Range("E20").Value = "  -1.42%  "
$ws.Range("D21").Value = "'12.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.58%  "
$ws.Range("D22").Value = "'69.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").Value = "'242.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.98%  "
$ws.Range("E24").Value = "  -2.94%  "
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "'25.39"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.63%  "
$ws.Range("E28").Value = "  -2.04%  "
$ws.Range("D29").Value = "'10.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.27%  "
$ws.Range("D30").Value = "'38.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.77%  "
$ws.Range("D31").Value = "'159.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "'5.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("E33").Value = "  +9.04%  "
$ws.Range("D34").Value = "'2.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "'0.0783"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("D36").Value = "'18.48"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("E37").Value = "  -7.56%  "
$ws.Range("E38").Value = "  -7.44%  "
$ws.Range("E39").Value = "  -1.52%  "
$ws.Range("D40").Value = "'0.117"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").Value = "'4.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").Value = "'22.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.15%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("D46").Value = "1.993.20"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").Value = "2.768.65"
$ws.Range("E48").Value = "  -1.67%  "
$ws.Range("E49").Value = "  -4.13%  "
$ws.Range("D50").Value = "'79.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.53%  "
$ws.Range("D51").Value = "'100.86"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.45%  "
